# Adds the new "Knärot - ekologi samt krav på livsmiljön" section (heading, body
# paragraphs and a "Referenser" subsection) at the very end of the document body -
# i.e. right after the last paragraph ("BILAGA 1 - Fridlysta arter") and before the
# trailing section properties - and bumps the date shown in the document's
# first-page header from 2023-09-13 to 2023-09-15.

$d = $word.ActiveDocument

# --- 1) Update the date shown in the document's first-page header. ---
# Headers(2) is the "first page" header (wdHeaderFooterFirstPage) which is the one
# that actually contains the "2023-09-13" / "Till:" / "Kopia: ..." text.

$section = $d.Sections(1)
$firstPageHeader = $section.Headers(2)
$firstPageHeader.Range.Find.Execute("2023-09-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-15", 2)

# --- 2) Append the new "Knärot" section at the end of the document body. ---
# The document body's last paragraph is the "BILAGA 1 - Fridlysta arter" title, so
# collapsing the whole-document range to its end and inserting OOXML there places
# the new content immediately after it (and before the final sectPr), matching the
# diff.

$insertionPoint = $d.Content
$insertionPoint.Collapse(0)

$knarotXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading1"/>
      </w:pPr>
      <w:r>
        <w:t>Knärot – ekologi samt krav på livsmiljön</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021).</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Samuel Johnsons doktorsavhandling </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
        </w:rPr>
        <w:t>“Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
        </w:rPr>
        <w:t xml:space="preserve">“Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">Vidare </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
        </w:rPr>
        <w:t>“More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II &amp; III).”</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
        </w:rPr>
        <w:t>“In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf &amp; Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022).</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022).</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading2"/>
      </w:pPr>
      <w:r>
        <w:t>Referenser - knärot</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">de Graaf M &amp; Roberts M.R., 2009. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
        </w:rPr>
        <w:t xml:space="preserve">Short-term response of the herbaceous layer within leave patches after harvest. </w:t>
      </w:r>
      <w:r>
        <w:t>Forest Ecology and Management 257, 1014-1025</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Halpern, C. B., Halaj, J., Evans, S. A., &amp; Dovciak, M., 2012. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
        </w:rPr>
        <w:t xml:space="preserve">Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">Ecological Applications, 22, 2049-2064 </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. &amp; Hylander, N. 2022. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
        </w:rPr>
        <w:t xml:space="preserve">Interactive effects of drought and edge exposure on old-growth forest understory species. </w:t>
      </w:r>
      <w:r>
        <w:t>Landscape Ecology, 37, sid 1839-1853</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Rudolphi, J., Jönsson, M. T., &amp; Gustafsson, L., 2014. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
        </w:rPr>
        <w:t xml:space="preserve">Biological legacies buffer local species extinction after logging. </w:t>
      </w:r>
      <w:r>
        <w:t>Journal of Applied Ecology. 51, 53-62.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Skogsstyrelsen, 2022. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
        </w:rPr>
        <w:t xml:space="preserve">Vägledning för hänsyn till knärot. </w:t>
      </w:r>
      <w:r>
        <w:t>https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">SLU Artdatabanken, 2021. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
        </w:rPr>
        <w:t xml:space="preserve">Artfaktablad. Naturvård – artfakta. </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">SLU Artdatabanken, Uppsala </w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertionPoint.InsertXML($knarotXml)
